$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Revert "By Section Enable/Disable": restore the original data grid A1:K15 ---

# 1) Write cell values (columns A..K = 1..11)
$ws.Cells.Item(1,2).Value = 10001
$ws.Cells.Item(1,3).Value = 10002
$ws.Cells.Item(1,4).Value = 10003
$ws.Cells.Item(1,5).Value = 10004
$ws.Cells.Item(1,6).Value = 10005
$ws.Cells.Item(1,7).Value = 10006
$ws.Cells.Item(1,8).Value = 10007
$ws.Cells.Item(1,9).Value = 10008
$ws.Cells.Item(1,10).Value = 10009
$ws.Cells.Item(1,11).Value = 10010
$ws.Cells.Item(2,1).Value = 10001
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(3,1).Value = 10002
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).Value = 0
$ws.Cells.Item(4,1).Value = 10003
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(5,1).Value = 10004
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = 0
$ws.Cells.Item(6,1).Value = 10005
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(6,8).Value = 0
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,10).Value = 0
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(7,1).Value = 10006
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 0
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = 0
$ws.Cells.Item(8,1).Value = 10007
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = 0
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = 0
$ws.Cells.Item(9,1).Value = 10008
$ws.Cells.Item(9,2).Value = 1
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,8).Value = 0
$ws.Cells.Item(9,9).Value = 1
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(10,1).Value = 10009
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 0
$ws.Cells.Item(10,9).Value = 0
$ws.Cells.Item(10,10).Value = 0
$ws.Cells.Item(10,11).Value = 0
$ws.Cells.Item(11,1).Value = 10010
$ws.Cells.Item(11,2).Value = 1
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 0
$ws.Cells.Item(11,9).Value = 1
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).Value = 1

# 2) Clear any leftover values beyond column K on the data rows (old layout used up to column G)
$ws.Range("L1:O15").ClearContents()

# 3) Re-apply the centered style (cellXfs index 1: horizontal=center, vertical=center)
#    to every cell in the restored range, matching the original workbook formatting.
$full = $ws.Range("A1:K15")
$full.HorizontalAlignment = -4108
$full.VerticalAlignment = -4108

# 4) A few rows in the original sheet also carried empty, centered cells out to column O
$ws.Range("L1:O1").HorizontalAlignment = -4108
$ws.Range("L1:O1").VerticalAlignment = -4108
$ws.Range("L7:O7").HorizontalAlignment = -4108
$ws.Range("L7:O7").VerticalAlignment = -4108
$ws.Range("L8:O8").HorizontalAlignment = -4108
$ws.Range("L8:O8").VerticalAlignment = -4108
$ws.Range("L10:O10").HorizontalAlignment = -4108
$ws.Range("L10:O10").VerticalAlignment = -4108

# 5) Restore the original active selection
$ws.Activate()
$ws.Range("B2").Select()
